# Fix for smoke test (Jenkins server)
# Append two new "Submission time" rows (10 and 11) to the
# "Checkertificate" sheet, extending the data range from A1:E9 to A1:E11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checkertificate")

# --- Row 10 ---
$ws.Cells.Item(10, 1).Value = "01.07.2023 00:15 (Kyiv+Israel) 22:15 (UTC) 07:15 (Japan) 03:45 (India)"
$ws.Cells.Item(10, 2).Value = "***"
$ws.Cells.Item(10, 3).Value = "***"
$ws.Cells.Item(10, 4).Value = 0.941
# Preserve the original floating point rounding noise (1.013 - 0.941)
$ws.Cells.Item(10, 5).Formula = "=1.013-0.941"
$ws.Cells.Item(10, 5).Value = $ws.Cells.Item(10, 5).Value()

# --- Row 11 ---
$ws.Cells.Item(11, 1).Value = "01.07.2023 00:21 (Kyiv+Israel) 22:21 (UTC) 07:21 (Japan) 03:51 (India)"
$ws.Cells.Item(11, 2).Value = "***"
$ws.Cells.Item(11, 3).Value = "***"
$ws.Cells.Item(11, 4).Value = 0.885
# Preserve the original floating point rounding noise (0.941 - 0.885)
$ws.Cells.Item(11, 5).Formula = "=0.941-0.885"
$ws.Cells.Item(11, 5).Value = $ws.Cells.Item(11, 5).Value()
